$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 06.25 (row 6) Goal: drop the "(완)" completion markers ---
$goal0625 = @"
1-1. Datebook Comment 부분 요약
1-2. GIT 보고서 다시 작성
1-3. PPT 머릿글 요약 작성
1-4. 파일명 수정
1-5 최종 아이디어 선정 이유 작성
2. 요구사항 수집 및 분석
"@
$ws.Range("D6").Value = $goal0625

# --- 06.28 (row 7): edit Time / Goal, add Comment ---
$time0628 = @"
08:30
18:00
"@
$ws.Range("C7").Value = $time0628
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("C7").VerticalAlignment = -4108
$ws.Range("C7").WrapText = $true

$goal0628 = @"
1. 06.28 회의록 작성
2. 요구사항 정의 및 분석
"@
$ws.Range("D7").Value = $goal0628

$comment0628 = @"
1. 진행사항 체크 및 2주차 목표 설정
2. 06.28 회의록 작성
3. 상용화 된 타 앱 기능 분석
4. 요구사항 정의
 - 회원가입, DB, 무결성 (미완)
5. vi사용법 설명
"@
$ws.Range("E7").Value = $comment0628
$ws.Range("E7").WrapText = $true

$ws.Rows(7).RowHeight = 79.05

# --- 06.29 (row 8): add a new Goal ---
$goal0629 = @"
1. vi 명령어 설명서 작성
2. 요구사항 정의 및 분석
"@
$ws.Range("D8").Value = $goal0629
$ws.Range("D8").WrapText = $true

$ws.Rows(8).RowHeight = 26.35

# --- selection moves to D7 ---
$ws.Range("D7").Select() | Out-Null
